$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LT Value Answer")

# Add a "PV Profit" check column (F) that mirrors column E's discounted
# profit formula, for rows 13-17 (the 5 forecast years).
$ws.Range("F13").Formula = "=D13/(1+`$A`$6)^A13"
$ws.Range("F14:F17").Formula = "=D14/(1+`$A`$6)^A14"

# Apply the same Currency accounting format used elsewhere in the table.
$curFmt = "_(`"$`"* #,##0.00_);_(`"$`"* \(#,##0.00\);_(`"$`"* `"-`"??_);_(@_)"
$ws.Range("F13:F17").NumberFormat = $curFmt
$ws.Range("F19").NumberFormat = $curFmt

# Match the author's last-saved cursor position.
$ws.Range("I16").Select()
